$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cashflow")

# Rescale the pasted-value cashflow figures (F2:BN9) by the revised factor.
for ($r = 2; $r -le 9; $r++) {
  for ($c = 6; $c -le 66; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.Value = $cell.Value() * -0.787421561
  }
}

# Remove the stray formatted-but-empty row that was left below the table.
$ws.Rows.Item(17).Delete()

# Restore the cursor to where the author left it.
[void]$ws.Activate()
[void]$ws.Range("A11").Select()
